$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is a "clean" numeric literal (e.g. "1.002",
# "0.08020") would otherwise be auto-coerced to a number by the Value
# setter, silently dropping significant trailing zeros / precision. Force
# them to Text format first so the assigned string is stored verbatim,
# matching the scraped source formatting exactly.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.253.46'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.805.06'
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '314.44'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D7").Value = '0.5252'
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("D8").Value = '0.3825'
$ws.Range("E8").Value = '  -2.36%  '
$ws.Range("D9").Value = '0.08020'
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.104'
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '41.43'
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").Value = '6.347'
$ws.Range("E12").Value = '  +1.55%  '
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").Value = '20.66'
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = '7.359'
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D16").Value = '1.802.20'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").Value = '92.73'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").Value = '0.00001099'
$ws.Range("E18").Value = '  -2.46%  '
$ws.Range("D19").Value = '0.06602'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").Value = '5.974'
$ws.Range("E22").Value = '  -1.88%  '
$ws.Range("D23").Value = '28.304.49'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = '2.233'
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("D26").Value = '160.31'
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").Value = '20.49'
$ws.Range("E27").Value = '  -2.56%  '
$ws.Range("D28").Value = '2.010.97'
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("D29").Value = '2.377'
$ws.Range("E29").Value = '  -0.86%  '
$ws.Range("D30").Value = '123.27'
$ws.Range("E30").Value = '  -1.86%  '
$ws.Range("D31").Value = '0.1083'
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = '1.061'
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").Value = '3.679'
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("D34").Value = '5.575'
$ws.Range("E34").Value = '  -1.67%  '
$ws.Range("D35").Value = '0.07251'
$ws.Range("E35").Value = '  +2.77%  '
$ws.Range("D36").Value = '12.42'
$ws.Range("E36").Value = '  +10.47%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.207'
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("D38").Value = '0.2168'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").Value = '0.02320'
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.852'
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").Value = '0.6225'
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("D42").Value = '1.165'
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").Value = '1.371'
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '13.30'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.6046'
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("D46").Value = '3.767'
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").Value = '127.19'
$ws.Range("E47").Value = '  +2.10%  '
$ws.Range("D48").Value = '1.217'
$ws.Range("E48").Value = '  +2.22%  '
$ws.Range("D49").Value = '1.936'
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").Value = '73.17'
$ws.Range("E51").Value = '  -1.64%  '

Write-Host "Applied cryptos update"
